# Add documentation rows for "virtualenv" commands to Sheet1 (rows 18-28),
# mirroring the existing COMMANDS / DESCRIPTION table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A/B pairs to append, starting at row 18 (existing data ends at row 17).
$rows = @(
    @("pip install virtualenv", "it will create a virtual environment"),
    @("virtualenv --version", "it show the version of virtualen"),
    @("virtualenv MyProject", "it will a virtualenv for MyProject"),
    @("cd/source/MyProject/Scripts/activate.bat", "it will activate the virualenv"),
    @("deactivate", "it will deactivate the virtual env"),
    @(" virtualenv --no-site-packages ", "it will not include the packages that are installed globally"),
    @("pip install virtualenvwrapper", "it will install virtualen wrapper"),
    @("mkvirtualenv my_project", "it will a virtualenv for MyProject"),
    @("workon my_project", "it will enable u to work on that project"),
    @("deactivate", "it will deactivate the virtual env"),
    @("rmvirtualenv venv", "it will remov/delete the virtualenv")
)

$startRow = 18
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $pair = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $pair[0]
    $ws.Cells.Item($r, 2).Value = $pair[1]
}

$lastRow = $startRow + $rows.Length - 1

# Update the view state to match: selection on the new last cell, and
# scroll the sheet down a bit so the new rows are in view.
$ws.Activate()
$ws.Cells.Item($lastRow, 2).Select()
$excel.ActiveWindow.ScrollRow = 14
